# Generate Report for handback
# Fills in the handoff/handback correspondence datetimes for the
# "fbdec5ff-..." file rows that were generated by the latest report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-09 04:27:15"
$wsZhCn.Range("G3").Value = "2016-01-09 04:27:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-09 04:27:25"
$wsDeDe.Range("G3").Value = "2016-01-09 04:28:15"
